$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Give the new last data row (30) the closing "last row" border/style ---
# (today row 31 -- the second worker's row, which is about to be removed -- carries that style;
# copy its formatting onto row 30 before deleting it)
$ws.Range("B31:J31").Copy()
$ws.Range("B30:J30").PasteSpecial(-4122)   # xlPasteFormats

# --- 2. Rewrite the period rows (16-30) for LEANDRO DIAZ ARRIETA in ascending period order ---
$periods = @("2108","2109","2110","2111","2112","2201","2202","2203","2204","2205","2206","2207","2208","2209","2210")
for ($i = 0; $i -lt $periods.Count; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    if ($i -eq $periods.Count - 1) {
        $ws.Cells.Item($row, 6).Value = 32707
    } else {
        $ws.Cells.Item($row, 6).Value = 36341
    }
}

# --- 3. Update the summary figures now that worker LUIS CARLOS JIMENEZ CASTILLO is gone ---
$ws.Range("E11").Value = 541481
$ws.Range("C13").Value = 1

# --- 4. Remove the second worker's row entirely (shifts the signature block up by one row) ---
$ws.Rows("31:31").Delete()

# --- 5. Column D ("Nombre Trabajador") no longer needs to fit the longer removed name ---
$ws.Columns("D:D").AutoFit()
